$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 5633.3335
$ws.Range("I21").Value = 6750
$ws.Range("J21").Value = 3400
$ws.Range("K21").Value = 6750
$ws.Range("L21").Value = 3400
$ws.Range("M21").Value = -6282
$ws.Range("N21").Value = -4336
$ws.Range("H23").Value = 5633.3335
$ws.Range("I23").Value = 6750
$ws.Range("J23").Value = 3400
$ws.Range("K23").Value = 6750
$ws.Range("L23").Value = 3400
$ws.Range("M23").Value = -6516
$ws.Range("N23").Value = -3868
$ws.Range("H43").Value = 5098.6
$ws.Range("I43").Value = 1997
$ws.Range("J43").Value = 7166.3335
$ws.Range("K43").Value = 1997
$ws.Range("L43").Value = 7166.3335
$ws.Range("M43").Value = -1928
$ws.Range("N43").Value = -7304.3335
$ws.Range("H51").Value = 11583.5
$ws.Range("J51").Value = 9900
$ws.Range("L51").Value = 9900
$ws.Range("N51").Value = -10868
$ws.Range("H53").Value = 117.117645
$ws.Range("I53").Value = 89.916664
$ws.Range("J53").Value = 182.4
$ws.Range("K53").Value = 89.916664
$ws.Range("L53").Value = 182.4
$ws.Range("M53").Value = 547.083336
$ws.Range("N53").Value = -1456.4
$ws.Range("H116").Value = 5847.25
$ws.Range("H141").Value = 4911.75
$ws.Range("I141").Value = 4549.1665
$ws.Range("K141").Value = 13647.4995
$ws.Range("M141").Value = -8467.499500000002

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 255
$ws.Range("I4").Value = 237.5
$ws.Range("J4").Value = 325
$ws.Range("K4").Value = 237.5
$ws.Range("L4").Value = 325
$ws.Range("M4").Value = -121.5
$ws.Range("N4").Value = -557
$ws.Range("H5").Value = 2000
$ws.Range("J5").Value = 2000
$ws.Range("L5").Value = 2000
$ws.Range("N5").Value = -2224
$ws.Range("H9").Value = 10000
$ws.Range("J9").Value = 10000
$ws.Range("L9").Value = 10000
$ws.Range("N9").Value = -10340
$ws.Range("H20").Value = 10000
$ws.Range("J20").Value = 10000
$ws.Range("L20").Value = 10000
$ws.Range("N20").Value = -10540
$ws.Range("H37").Value = 14749.5
$ws.Range("J37").Value = 21999
$ws.Range("L37").Value = 21999
$ws.Range("N37").Value = -22545
$ws.Range("H44").Value = 39997
$ws.Range("J44").Value = 39997
$ws.Range("L44").Value = 39997
$ws.Range("N44").Value = -40973
$ws.Range("H45").Value = 2881.1538
$ws.Range("I45").Value = 1924.5714
$ws.Range("K45").Value = 1924.5714
$ws.Range("M45").Value = -1547.5714
$ws.Range("H55").Value = 25010
$ws.Range("J55").Value = 31997.334
$ws.Range("L55").Value = 31997.334
$ws.Range("N55").Value = -32627.334
$ws.Range("H61").Value = 6213.357
$ws.Range("J61").Value = 4664.6665
$ws.Range("L61").Value = 4664.6665
$ws.Range("N61").Value = -5088.6665
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H97").Value = 602.6429000000001
$ws.Range("I97").Value = 582.25
$ws.Range("J97").Value = 725
$ws.Range("K97").Value = 582.25
$ws.Range("L97").Value = 725
$ws.Range("M97").Value = -86.25
$ws.Range("N97").Value = -1717
$ws.Range("H122").Value = 2054.5
$ws.Range("I122").Value = 2054.5
$ws.Range("K122").Value = 6163.5
$ws.Range("M122").Value = -3713.5
$ws.Range("H132").Value = 1930.8158
$ws.Range("I132").Value = 1929.9375
$ws.Range("K132").Value = 5789.8125
$ws.Range("M132").Value = -3259.8125
$ws.Range("H136").Value = 6213.357
$ws.Range("J136").Value = 4664.6665
$ws.Range("L136").Value = 13993.9995
$ws.Range("N136").Value = -19093.9995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 2000
$ws.Range("J4").Value = 2000
$ws.Range("L4").Value = 2000
$ws.Range("N4").Value = -2230
$ws.Range("H86").Value = 942.44446
$ws.Range("I86").Value = 1034.625
$ws.Range("J86").Value = 868.7
$ws.Range("K86").Value = 1034.625
$ws.Range("L86").Value = 868.7
$ws.Range("M86").Value = 88.375
$ws.Range("N86").Value = -3114.7
$ws.Range("H89").Value = 942.44446
$ws.Range("I89").Value = 1034.625
$ws.Range("J89").Value = 868.7
$ws.Range("K89").Value = 5173.125
$ws.Range("L89").Value = 4343.5
$ws.Range("M89").Value = 442.875
$ws.Range("N89").Value = -15575.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3555.25
$ws.Range("J16").Value = 4500
$ws.Range("L16").Value = 4500
$ws.Range("N16").Value = -5074
$ws.Range("H22").Value = 59190.145
$ws.Range("J22").Value = 4907.231
$ws.Range("L22").Value = 4907.231
$ws.Range("N22").Value = -5607.231
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H105").Value = 622.8889
$ws.Range("I105").Value = 638.25
$ws.Range("K105").Value = 638.25
$ws.Range("M105").Value = 1108.75
$ws.Range("H113").Value = 3555.25
$ws.Range("J113").Value = 4500
$ws.Range("L113").Value = 4500
$ws.Range("N113").Value = -8840

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 161.72728
$ws.Range("I5").Value = 161.72728
$ws.Range("K5").Value = 485.18184
$ws.Range("M5").Value = -373.18184
$ws.Range("H60").Value = 50
$ws.Range("I60").Value = 50
$ws.Range("K60").Value = 150
$ws.Range("M60").Value = 101
$ws.Range("H129").Value = 825.7143
$ws.Range("I129").Value = 776
$ws.Range("K129").Value = 2328
$ws.Range("M129").Value = 2672
$ws.Range("H135").Value = 161.72728
$ws.Range("I135").Value = 161.72728
$ws.Range("K135").Value = 1455.54552
$ws.Range("M135").Value = 1079.45448

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 3425.5
$ws.Range("I5").Value = 3425.5
$ws.Range("K5").Value = 3425.5
$ws.Range("M5").Value = -3313.5
$ws.Range("H11").Value = 117333390
$ws.Range("I11").Value = 117333390
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 117333390
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -117333251
$ws.Range("N11").ClearContents()
$ws.Range("H35").Value = 2947480
$ws.Range("J35").Value = 1750000
$ws.Range("L35").Value = 1750000
$ws.Range("N35").Value = -1750596
$ws.Range("H70").Value = 16041.5
$ws.Range("I70").Value = 12499.429
$ws.Range("K70").Value = 12499.429
$ws.Range("M70").Value = -12229.429
$ws.Range("H73").Value = 16041.5
$ws.Range("I73").Value = 12499.429
$ws.Range("K73").Value = 12499.429
$ws.Range("M73").Value = -11563.429

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 24855.572
$ws.Range("I47").Value = 24000
$ws.Range("J47").Value = 24998.166
$ws.Range("K47").Value = 24000
$ws.Range("L47").Value = 24998.166
$ws.Range("M47").Value = -23510
$ws.Range("N47").Value = -25978.166
$ws.Range("H52").Value = 24855.572
$ws.Range("I52").Value = 24000
$ws.Range("J52").Value = 24998.166
$ws.Range("K52").Value = 24000
$ws.Range("L52").Value = 24998.166
$ws.Range("M52").Value = -23767
$ws.Range("N52").Value = -25464.166
$ws.Range("H55").Value = 688.3077
$ws.Range("I55").Value = 720.875
$ws.Range("K55").Value = 720.875
$ws.Range("M55").Value = -547.875

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 400500
$ws.Range("J2").Value = 625
$ws.Range("L2").Value = 625
$ws.Range("N2").Value = -849
